$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.03900598280555201
$ws.Range("D2").Value = 0.06171644822111944
$ws.Range("E2").Value = 0.1273645648577855
$ws.Range("F2").Value = 3.024660653864856
$ws.Range("G2").Value = 0.002548312747561555
$ws.Range("I2").Value = 1.795120758473388
$ws.Range("J2").Value = 0.2213205928045738
$ws.Range("K2").Value = 2.949830790324086
$ws.Range("N2").Value = 1.68198484549966
$ws.Range("C3").Value = 0.03908443852522403
$ws.Range("D3").Value = 0.06036937914596052
$ws.Range("E3").Value = 0.1247593813203594
$ws.Range("F3").Value = 3.002578464191458
$ws.Range("G3").Value = 0.002554012195755754
$ws.Range("I3").Value = 1.779397174936918
$ws.Range("J3").Value = 0.216842663057264
$ws.Range("K3").Value = 2.772341838323598
$ws.Range("N3").Value = 1.705569986164299
$ws.Range("C4").Value = 0.0391360612526519
$ws.Range("D4").Value = 0.05955873523247845
$ws.Range("E4").Value = 0.1232267231857165
$ws.Range("F4").Value = 2.990965839343943
$ws.Range("G4").Value = 0.002557693423782093
$ws.Range("I4").Value = 1.770961239297279
$ws.Range("J4").Value = 0.2142320986826718
$ws.Range("K4").Value = 2.664938426658011
$ws.Range("N4").Value = 1.720737386578969
$ws.Range("C5").Value = 0.03915797147954514
$ws.Range("D5").Value = 0.05923256382047981
$ws.Range("E5").Value = 0.1226189460210279
$ws.Range("F5").Value = 2.98672073641238
$ws.Range("G5").Value = 0.002559239418056618
$ws.Range("I5").Value = 1.76782834833574
$ws.Range("J5").Value = 0.2132030283997324
$ws.Range("K5").Value = 2.621564533163905
$ws.Range("N5").Value = 1.727090438983241
$ws.Range("C6").Value = 0.03916166259216247
$ws.Range("D6").Value = 0.05917865630254937
$ws.Range("E6").Value = 0.1225190382935644
$ws.Range("F6").Value = 2.986045198814608
$ws.Range("G6").Value = 0.002559498903807898
$ws.Range("I6").Value = 1.767326499446895
$ws.Range("J6").Value = 0.2130342462550132
$ws.Range("K6").Value = 2.614386053837336
$ws.Range("N6").Value = 1.728155751308082
$ws.Range("C7").Value = 0.03913635319787367
$ws.Range("D7").Value = 0.05955431943159084
$ws.Range("E7").Value = 0.1232184585422402
$ws.Range("F7").Value = 2.990906618916341
$ws.Range("G7").Value = 0.002557714087539959
$ws.Range("I7").Value = 1.770917755897699
$ws.Range("J7").Value = 0.2142180797801743
$ws.Range("K7").Value = 2.664351879425396
$ws.Range("N7").Value = 1.720822369187355
$ws.Range("C8").Value = 0.03903232163113657
$ws.Range("D8").Value = 0.0612485812914656
$ws.Range("E8").Value = 0.1264523808795772
$ws.Range("F8").Value = 3.016641500145198
$ws.Range("G8").Value = 0.002550240294792983
$ws.Range("I8").Value = 1.789445434081912
$ws.Range("J8").Value = 0.2197476670035883
$ws.Range("K8").Value = 2.88830413174685
$ws.Range("N8").Value = 1.689974440965079
$ws.Range("C9").Value = 0.03885544414278286
$ws.Range("D9").Value = 0.06470041834172235
$ws.Range("E9").Value = 0.133327442022015
$ws.Range("F9").Value = 3.082653752521992
$ws.Range("G9").Value = 0.00253701880708087
$ws.Range("I9").Value = 1.835520856731875
$ws.Range("J9").Value = 0.2317020068855697
$ws.Range("K9").Value = 3.34010703526269
$ws.Range("N9").Value = 1.634937371502317
$ws.Range("C10").Value = 0.03874170301394386
$ws.Range("D10").Value = 0.06731406027853382
$ws.Range("E10").Value = 0.1387075901179884
$ws.Range("F10").Value = 3.140788450616668
$ws.Range("G10").Value = 0.002528169081794911
$ws.Range("I10").Value = 1.875423131676769
$ws.Range("J10").Value = 0.2411754659610068
$ws.Range("K10").Value = 3.679971656325563
$ws.Range("N10").Value = 1.59784317808376
$ws.Range("C11").Value = 0.03869341368410417
$ws.Range("D11").Value = 0.06851967292529793
$ws.Range("E11").Value = 0.141227520282051
$ws.Range("F11").Value = 3.169363032387054
$ws.Range("G11").Value = 0.002524328504911386
$ws.Range("I11").Value = 1.894914348718942
$ws.Range("J11").Value = 0.2456381724753953
$ws.Range("K11").Value = 3.836357481203038
$ws.Range("N11").Value = 1.581697138953269
$ws.Range("C12").Value = 0.03867561915946283
$ws.Range("D12").Value = 0.06897857559913234
$ws.Range("E12").Value = 0.1421922359687144
$ws.Range("F12").Value = 3.180492351292202
$ws.Range("G12").Value = 0.002522900639322468
$ws.Range("I12").Value = 1.902489713290962
$ws.Range("J12").Value = 0.2473503370941046
$ws.Range("K12").Value = 3.895836158178497
$ws.Range("N12").Value = 1.575688139996231
$ws.Range("C13").Value = 0.03867942973944594
$ws.Range("D13").Value = 0.06887963808290465
$ws.Range("E13").Value = 0.1419840005689608
$ws.Range("F13").Value = 3.178081679558147
$ws.Range("G13").Value = 0.00252320698055285
$ws.Range("I13").Value = 1.900849544399478
$ws.Range("J13").Value = 0.2469805995828267
$ws.Range("K13").Value = 3.883014806355391
$ws.Range("N13").Value = 1.576977601672965
$ws.Range("C14").Value = 0.03869193989307718
$ws.Range("D14").Value = 0.06855737993871003
$ws.Range("E14").Value = 0.1413066778314587
$ws.Range("F14").Value = 3.170272445257837
$ws.Range("G14").Value = 0.002524210503789486
$ws.Range("I14").Value = 1.895533670484184
$ws.Range("J14").Value = 0.2457785864754385
$ws.Range("K14").Value = 3.841245625262445
$ws.Range("N14").Value = 1.581200665028511
$ws.Range("C15").Value = 0.03869966660518287
$ws.Range("D15").Value = 0.06836029437184266
$ws.Range("E15").Value = 0.1408931635184416
$ws.Range("F15").Value = 3.165529350090424
$ws.Range("G15").Value = 0.002524828634017532
$ws.Range("I15").Value = 1.892302924510474
$ws.Range("J15").Value = 0.2450452205453502
$ws.Range("K15").Value = 3.815694579614444
$ws.Range("N15").Value = 1.583801120524049
$ws.Range("C16").Value = 0.03874492788411033
$ws.Range("D16").Value = 0.0672356030698964
$ws.Range("E16").Value = 0.1385443694471817
$ws.Range("F16").Value = 3.138964094997448
$ws.Range("G16").Value = 0.002528423786039194
$ws.Range("I16").Value = 1.874176445820865
$ws.Range("J16").Value = 0.2408869176252182
$ws.Range("K16").Value = 3.669787527877588
$ws.Range("N16").Value = 1.598913051174506
$ws.Range("C17").Value = 0.0387735748864344
$ws.Range("D17").Value = 0.06654988422874908
$ws.Range("E17").Value = 0.1371220584840245
$ws.Range("F17").Value = 3.123214265205291
$ws.Range("G17").Value = 0.002530676620945098
$ws.Range("I17").Value = 1.863400815769239
$ws.Range("J17").Value = 0.2383753132248785
$ws.Range("K17").Value = 3.580736095963232
$ws.Range("N17").Value = 1.608370565740977
$ws.Range("C18").Value = 0.03879037708490074
$ws.Range("D18").Value = 0.06615704626595686
$ws.Range("E18").Value = 0.1363108012441288
$ws.Range("F18").Value = 3.114355594405652
$ws.Range("G18").Value = 0.002531989833900315
$ws.Range("I18").Value = 1.857328923516377
$ws.Range("J18").Value = 0.2369451104362099
$ws.Range("K18").Value = 3.529683406767674
$ws.Range("N18").Value = 1.613878802346218
$ws.Range("C19").Value = 0.03879612202434402
$ws.Range("D19").Value = 0.06602430859283004
$ws.Range("E19").Value = 0.1360372930204861
$ws.Range("F19").Value = 3.111390508502041
$ws.Range("G19").Value = 0.002532437465460472
$ws.Range("I19").Value = 1.855294663873167
$ws.Range("J19").Value = 0.2364633362842881
$ws.Range("K19").Value = 3.512426498932371
$ws.Range("N19").Value = 1.615755552882186
$ws.Range("C20").Value = 0.03877049174333536
$ws.Range("D20").Value = 0.06662271797645047
$ws.Range("E20").Value = 0.1372727598699868
$ws.Range("F20").Value = 3.124870121150423
$ws.Range("G20").Value = 0.002530434998873596
$ws.Range("I20").Value = 1.864534849770152
$ws.Range("J20").Value = 0.2386411853856316
$ws.Range("K20").Value = 3.590198433654052
$ws.Range("N20").Value = 1.607356702099841
$ws.Range("C21").Value = 0.03868825205295678
$ws.Range("D21").Value = 0.06865197107508436
$ws.Range("E21").Value = 0.1415053392918111
$ws.Range("F21").Value = 3.172557806288324
$ws.Range("G21").Value = 0.00252391502734266
$ws.Range("I21").Value = 1.897089778938749
$ws.Range("J21").Value = 0.2461310417533582
$ws.Range("K21").Value = 3.853507200539639
$ws.Range("N21").Value = 1.57995739159208
$ws.Range("C22").Value = 0.03863736696662201
$ws.Range("D22").Value = 0.06999196643571537
$ws.Range("E22").Value = 0.1443326433128291
$ws.Range("F22").Value = 3.205525316538626
$ws.Range("G22").Value = 0.002519808108932832
$ws.Range("I22").Value = 1.919500669472569
$ws.Range("J22").Value = 0.25115579466447
$ws.Range("K22").Value = 4.027104327009795
$ws.Range("N22").Value = 1.562663492417053
$ws.Range("C23").Value = 0.0386642648379123
$ws.Range("D23").Value = 0.06927553677405029
$ws.Range("E23").Value = 0.1428180527700462
$ws.Range("F23").Value = 3.187764279275711
$ws.Range("G23").Value = 0.002521985984895897
$ws.Range("I23").Value = 1.907435142509229
$ws.Range("J23").Value = 0.2484620525358849
$ws.Range("K23").Value = 3.934313181765901
$ws.Range("N23").Value = 1.571837325777061
$ws.Range("C24").Value = 0.03877188459570391
$ws.Range("D24").Value = 0.06658978550400008
$ws.Range("E24").Value = 0.1372046076923468
$ws.Range("F24").Value = 3.124120898017509
$ws.Range("G24").Value = 0.002530544180097071
$ws.Range("I24").Value = 1.864021769518743
$ws.Range("J24").Value = 0.2385209417626442
$ws.Range("K24").Value = 3.585920061954312
$ws.Range("N24").Value = 1.607814848780388
$ws.Range("C25").Value = 0.03890042453580378
$ws.Range("D25").Value = 0.06375288967567627
$ws.Range("E25").Value = 0.1314100605271342
$ws.Range("F25").Value = 3.063115165124117
$ws.Range("G25").Value = 0.002540443070317781
$ws.Range("I25").Value = 1.82200180900648
$ws.Range("J25").Value = 0.2283477460436671
$ws.Range("K25").Value = 3.216507099837315
$ws.Range("N25").Value = 1.64924071885727
